$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "Yu Qiao"
$ws.Range("C1").Value = "Luc Van Gool"
$ws.Range("D1").Value = "Lei Zhang"

$ws.Range("A5").Value = "Total"
